# Regenerate merged AHB files:
#  - rename the "_old"/"_new" comparison-column headers to the concrete
#    format versions being diffed (FV2210 -> FV2304)
#  - freeze the header row
#  - wrap the data range in a native Excel table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A-J (1-10) held the "*_old" headers -> rename to "*_FV2210"
for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# Column K (11) stays "diff"; columns L-U (12-21) held "*_new" -> "*_FV2304"
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# Freeze the header row (row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a native table so it gets autofilter + styling
$tableRange = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
